$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "Procent"

# Formula for D2:D13 -> percent of total (107)
# (Entered as two assignments so D2 stays a standalone formula while
# D3:D13 form a single shared-formula group, matching the authored file.)
$ws.Range("D2").Formula = "=C2/107*100"
$ws.Range("D3:D13").Formula = "=C3/107*100"

# Update selection to match target (F10)
$ws.Range("F10").Select()
